# Auto-generated PowerShell-style Excel COM-interop script
# Updates the "Price" (D) and "Volume(1h)" (E) columns of the cryptos sheet
# to reflect refreshed values from the Mon May 15 03:29:01 UTC 2023 data pull.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.719.45"
$ws.Range("E2").Value = "  +0.62%  "
$ws.Range("D3").Value = "1.863.30"
$ws.Range("E3").Value = "  +0.72%  "
$ws.Range("E4").Value = "  -0.61%  "
$ws.Range("D5").Value = "'320.77"
$ws.Range("E5").Value = "  -0.26%  "
$ws.Range("E6").Value = "  -0.85%  "
$ws.Range("D7").Value = "'0.4360"
$ws.Range("E7").Value = "  -0.57%  "
$ws.Range("D8").Value = "'0.3799"
$ws.Range("E8").Value = "  +0.71%  "
$ws.Range("D9").Value = "'0.07435"
$ws.Range("E9").Value = "  +0.40%  "
$ws.Range("D10").Value = "'0.8844"
$ws.Range("E10").Value = "  +1.09%  "
$ws.Range("D11").Value = "'21.63"
$ws.Range("E11").Value = "  +0.71%  "
$ws.Range("D12").Value = "1.859.79"
$ws.Range("E12").Value = "  +0.26%  "
$ws.Range("D13").Value = "'6.751"
$ws.Range("E13").Value = "  +0.86%  "
$ws.Range("E14").Value = "  -0.59%  "
$ws.Range("E15").Value = "  -1.04%  "
$ws.Range("D16").Value = "'86.82"
$ws.Range("E16").Value = "  +4.82%  "
$ws.Range("D17").Value = "'1.025"
$ws.Range("E17").Value = "  -1.04%  "
$ws.Range("D18").Value = "'0.000009076"
$ws.Range("E18").Value = "  +0.59%  "
$ws.Range("D19").Value = "'1.020"
$ws.Range("E19").Value = "  -0.84%  "
$ws.Range("D20").Value = "'15.47"
$ws.Range("E20").Value = "  +0.39%  "
$ws.Range("D21").Value = "27.715.29"
$ws.Range("E21").Value = "  +0.50%  "
$ws.Range("D22").Value = "'5.292"
$ws.Range("E22").Value = "  +0.74%  "
$ws.Range("E23").Value = "  -1.56%  "
$ws.Range("D24").Value = "2.100.78"
$ws.Range("E24").Value = "  +1.21%  "
$ws.Range("D25").Value = "'2.040"
$ws.Range("E25").Value = "  +6.25%  "
$ws.Range("D26").Value = "'157.43"
$ws.Range("E26").Value = "  -0.19%  "
$ws.Range("D27").Value = "'18.74"
$ws.Range("E27").Value = "  +0.13%  "
$ws.Range("D28").Value = "'5.372"
$ws.Range("E28").Value = "  +2.27%  "
$ws.Range("D29").Value = "'1.993"
$ws.Range("E29").Value = "  +1.30%  "
$ws.Range("D30").Value = "'120.73"
$ws.Range("E30").Value = "  +3.26%  "
$ws.Range("D31").Value = "'0.09055"
$ws.Range("E31").Value = "  +0.08%  "
$ws.Range("D32").Value = "'1.221"
$ws.Range("E32").Value = "  +2.28%  "
$ws.Range("D33").Value = "'0.7695"
$ws.Range("E33").Value = "  +1.16%  "
$ws.Range("D34").Value = "'3.032"
$ws.Range("E34").Value = "  +5.16%  "
$ws.Range("E35").Value = "  +1.42%  "
$ws.Range("E36").Value = "  -0.68%  "
$ws.Range("D37").Value = "'1.143"
$ws.Range("E37").Value = "  -0.41%  "
$ws.Range("D38").Value = "'0.01979"
$ws.Range("E38").Value = "  +0.36%  "
$ws.Range("D39").Value = "'0.05302"
$ws.Range("E39").Value = "  +0.16%  "
$ws.Range("D40").Value = "'2.880"
$ws.Range("E40").Value = "  +3.08%  "
$ws.Range("E41").Value = "  +0.77%  "
$ws.Range("D42").Value = "'6.962"
$ws.Range("E42").Value = "  +3.54%  "
$ws.Range("D43").Value = "'0.1681"
$ws.Range("E43").Value = "  +0.42%  "
$ws.Range("D44").Value = "'8.713"
$ws.Range("E44").Value = "  +2.85%  "
$ws.Range("D45").Value = "'10.80"
$ws.Range("E45").Value = "  +2.83%  "
$ws.Range("D46").Value = "'110.10"
$ws.Range("E46").Value = "  +1.31%  "
$ws.Range("D47").Value = "'1.714"
$ws.Range("E47").Value = "  +0.47%  "
$ws.Range("E48").Value = "  -0.96%  "
$ws.Range("D49").Value = "'0.06511"
$ws.Range("E49").Value = "  +2.03%  "
$ws.Range("D50").Value = "'0.4725"
$ws.Range("E50").Value = "  +1.68%  "
$ws.Range("E51").Value = "  +0.94%  "
